$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 27) with the latest portfolio figures.
# Force column A to be treated as text so the date-like string "2025-09-11"
# is stored as a plain string (matching the existing rows) instead of being
# auto-converted into a date serial number by Excel.
$row = 27
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "2025-09-11"
$ws.Range("A$row").Style = "Normal"

$ws.Range("B$row").Value = 56.93999862670898
$ws.Range("C$row").Value = 705.8499755859375
$ws.Range("D$row").Value = 328.1000061035156
